# Apply updated cryptocurrency price/volume data to sheet1 (ActiveSheet)
# The diff updates text-valued Price (D) and Volume(1h) (E) cells for rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the target range to Text format first so Excel does not
# auto-convert numeric-looking strings (e.g. "589.05") into real numbers,
# then clear the formatting afterwards so the cell style matches the original
# (plain, unstyled) cells -- only the stored text content should change.
$updateRange = $ws.Range("D2:E51")
$updateRange.NumberFormat = "@"

$ws.Range("D2").Value = "63.989.58"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").Value = "3.133.76"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "589.05"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").Value = "147.17"
$ws.Range("E6").Value = "  +2.34%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.127.41"
$ws.Range("E8").Value = "  +0.75%  "
$ws.Range("D9").Value = "0.533"
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("E10").Value = "  +12.37%  "
$ws.Range("D11").Value = "5.76"
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("E13").Value = "  +4.24%  "
$ws.Range("D14").Value = "37.38"
$ws.Range("E14").Value = "  +4.96%  "
$ws.Range("E15").Value = "  -1.06%  "
$ws.Range("D16").Value = "3.652.67"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").Value = "63.835.63"
$ws.Range("E17").Value = "  +1.55%  "
$ws.Range("D18").Value = "7.17"
$ws.Range("E18").Value = "  -1.98%  "
$ws.Range("D19").Value = "3.134.00"
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("D20").Value = "466.69"
$ws.Range("E20").Value = "  +2.54%  "
$ws.Range("E21").Value = "  +1.60%  "
$ws.Range("D22").Value = "0.735"
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("E23").Value = "  +0.75%  "
$ws.Range("D24").Value = "13.27"
$ws.Range("E24").Value = "  -3.75%  "
$ws.Range("D25").Value = "82.53"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  +8.49%  "
$ws.Range("E28").Value = "  +0.65%  "
$ws.Range("E29").Value = "  -1.88%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").Value = "6.87"
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("D32").Value = "27.16"
$ws.Range("E32").Value = "  -0.06%  "
$ws.Range("E33").Value = "  -3.50%  "
$ws.Range("D34").Value = "0.0₃0886"
$ws.Range("E34").Value = "  +10.26%  "
$ws.Range("E35").Value = "  +7.98%  "
$ws.Range("E36").Value = "  +0.95%  "
$ws.Range("D37").Value = "3.43"
$ws.Range("E37").Value = "  +11.73%  "
$ws.Range("D38").Value = "6.12"
$ws.Range("E38").Value = "  +0.40%  "
$ws.Range("D39").Value = "458.58"
$ws.Range("E39").Value = "  +7.24%  "
$ws.Range("D40").Value = "50.94"
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("D41").Value = "8.72"
$ws.Range("E41").Value = "  -1.55%  "
$ws.Range("D42").Value = "0.0373"
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").Value = "2.901.25"
$ws.Range("E43").Value = "  -1.17%  "
$ws.Range("D44").Value = "0.278"
$ws.Range("E44").Value = "  -0.43%  "
$ws.Range("E45").Value = "  +1.18%  "
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").Value = "35.85"
$ws.Range("E47").Value = "  +0.92%  "
$ws.Range("D48").Value = "126.03"
$ws.Range("E48").Value = "  +1.85%  "
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("D51").Value = "24.78"
$ws.Range("E51").Value = "  +0.16%  "

$updateRange.ClearFormats()
